# Update Sheet1: set G3, G4, G6 to 0 (formulas in H column recompute automatically)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("G3").Value = 0
$ws1.Range("G4").Value = 0
$ws1.Range("G6").Value = 0

# Update selection on Sheet1 to G7
$ws1.Range("G7").Select()

# Add a new worksheet named "alt" after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "alt"

# Row 1 headers
$ws2.Range("B1").Value = "pct of total"
$ws2.Range("C1").Value = "points"
$ws2.Range("D1").Value = "score"
$ws2.Range("E1").Value = "contrib to total"
$ws2.Range("F1").Value = "score"
$ws2.Range("G1").Value = "contrib to total"
$ws2.Range("H1").Value = "score"
$ws2.Range("I1").Value = "contrib to total"

$ws2.Range("D1:E1").Style = "20% - Accent1"
$ws2.Range("F1:G1").Style = "40% - Accent1"
$ws2.Range("H1:I1").Style = "20% - Accent1"

# Row 2 - P/A
$ws2.Range("A2").Value = "P/A"
$ws2.Range("B2").Value = 0.1
$ws2.Range("C2").Formula = "=100*B2"
$ws2.Range("D2").Value = 95
$ws2.Range("E2").Formula = "=(D2+5)*B2"
$ws2.Range("F2").Value = 55
$ws2.Range("G2").Formula = "=(F2+5)*`$B`$2"
$ws2.Range("H2").Value = 55
$ws2.Range("I2").Formula = "=H2*`$B`$2"

$ws2.Range("D2:E2").Style = "20% - Accent1"
$ws2.Range("F2:G2").Style = "40% - Accent1"
$ws2.Range("H2:I2").Style = "20% - Accent1"

# Rows 3-10 OA1-OA5, J, E best, E
$names = @("OA1", "OA2", "OA3", "OA4", "OA5", "J", "E best", "E ")
$bvals = @(0.04, 0.04, 0.04, 0.04, 0.04, 0.05, 0.35, 0.3)

for ($i = 0; $i -lt 8; $i++) {
    $r = 3 + $i
    $ws2.Range("A$r").Value = $names[$i]
    $ws2.Range("B$r").Value = $bvals[$i]
    $ws2.Range("C$r").Formula = "=100*B$r"
    $ws2.Range("D$r").Value = 95
    $ws2.Range("E$r").Formula = "=(D$r+5)*B$r"
    $ws2.Range("F$r").Value = 55
    $ws2.Range("G$r").Formula = "=F$r*`$B`$$r"
    $ws2.Range("H$r").Value = 55
    $ws2.Range("I$r").Formula = "=H$r*`$B`$$r"

    $ws2.Range("D$r`:E$r").Style = "20% - Accent1"
    $ws2.Range("F$r`:G$r").Style = "40% - Accent1"
    $ws2.Range("H$r`:I$r").Style = "20% - Accent1"
}

# Row 11 - total
$ws2.Range("A11").Value = "total"
$ws2.Range("B11").Formula = "=SUM(B2:B10)"
$ws2.Range("C11").Formula = "=SUM(C2:C10)"
$ws2.Range("E11").Formula = "=SUM(E2:E10)"
$ws2.Range("G11").Formula = "=SUM(G2:G10)"
$ws2.Range("I11").Formula = "=SUM(I2:I10)"

$ws2.Range("D11:E11").Style = "20% - Accent1"
$ws2.Range("F11:G11").Style = "40% - Accent1"
$ws2.Range("H11:I11").Style = "20% - Accent1"

# Row 16-17
$ws2.Range("B16").Value = "each oa"
$ws2.Range("C16").Value = 4
$ws2.Range("B17").Value = "each journal"
$ws2.Range("C17").Formula = "=5/15"

$ws2.Range("E12").Select()

$ws1.Select()
$ws1.Range("G7").Select()
